$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 7013179402
$ws.Range("D2").Value = 797738218

$ws.Range("C3").Value = 2894233510
$ws.Range("D3").Value = 318676377

$ws.Range("C4").Value = 1926715299
$ws.Range("D4").Value = 204210511

$ws.Range("C5").Value = 1927371690
$ws.Range("D5").Value = 196916027

$ws.Range("C6").Value = 3776801802
$ws.Range("D6").Value = 377881831

$ws.Range("C7").Value = 7133490408
$ws.Range("D7").Value = 693993467

$ws.Range("C8").Value = 6977967972
$ws.Range("D8").Value = 664489461

$ws.Range("E2:E8").Select()
